# Generate Report for Handoff
# Updates the "Latest HO Xliff Generate Date" (Overview sheet) and the
# "Latest Handoff Datetime" (per-locale sheets) for the file
# 8a44410b-c6f6-49fc-a00f-c20408a33d1e.md now that a new handoff xliff
# has been generated/handed off.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# Overview sheet: row 5 corresponds to 8a44410b-c6f6-49fc-a00f-c20408a33d1e.md
$wsOverview.Range("G5").Value = "2016-11-15 16:26:26"

# zh-cn sheet: row 5 corresponds to 8a44410b-c6f6-49fc-a00f-c20408a33d1e.md
$wsZhCn.Range("H5").Value = "2016-11-15 16:26:13"

# de-de sheet: row 5 corresponds to 8a44410b-c6f6-49fc-a00f-c20408a33d1e.md
$wsDeDe.Range("H5").Value = "2016-11-15 16:26:26"

$wb.Save()
